$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions) - update column F (想去人数 / interest count)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 819
$ws1.Cells.Item(4, 6).Value = 306
$ws1.Cells.Item(7, 6).Value = 332
$ws1.Cells.Item(8, 6).Value = 52
$ws1.Cells.Item(9, 6).Value = 124
$ws1.Cells.Item(10, 6).Value = 127
$ws1.Cells.Item(11, 6).Value = 1202
$ws1.Cells.Item(14, 6).Value = 884
$ws1.Cells.Item(15, 6).Value = 876
$ws1.Cells.Item(17, 6).Value = 70
$ws1.Cells.Item(18, 6).Value = 74
$ws1.Cells.Item(20, 6).Value = 772
$ws1.Cells.Item(21, 6).Value = 1741
$ws1.Cells.Item(22, 6).Value = 2863
$ws1.Cells.Item(23, 6).Value = 833
$ws1.Cells.Item(25, 6).Value = 2199
$ws1.Cells.Item(26, 6).Value = 670
$ws1.Cells.Item(27, 6).Value = 3040
$ws1.Cells.Item(28, 6).Value = 590
$ws1.Cells.Item(29, 6).Value = 13
$ws1.Cells.Item(30, 6).Value = 12
$ws1.Cells.Item(32, 6).Value = 734
$ws1.Cells.Item(34, 6).Value = 131
$ws1.Cells.Item(36, 6).Value = 1072
$ws1.Cells.Item(37, 6).Value = 1769
$ws1.Cells.Item(38, 6).Value = 389
$ws1.Cells.Item(40, 6).Value = 554
$ws1.Cells.Item(41, 6).Value = 187
$ws1.Cells.Item(42, 6).Value = 132
$ws1.Cells.Item(43, 6).Value = 173
$ws1.Cells.Item(44, 6).Value = 43

# Sheet 2: 演出 (Performances) - update column F
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 139
$ws2.Cells.Item(9, 6).Value = 13
$ws2.Cells.Item(12, 6).Value = 78

# Sheet 4: 全部类型 (All types) - update column F
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 819
$ws4.Cells.Item(4, 6).Value = 306
$ws4.Cells.Item(7, 6).Value = 332
$ws4.Cells.Item(8, 6).Value = 52
$ws4.Cells.Item(9, 6).Value = 124
$ws4.Cells.Item(10, 6).Value = 1202
$ws4.Cells.Item(12, 6).Value = 884
$ws4.Cells.Item(13, 6).Value = 876
$ws4.Cells.Item(14, 6).Value = 139
$ws4.Cells.Item(16, 6).Value = 70
$ws4.Cells.Item(18, 6).Value = 74
$ws4.Cells.Item(19, 6).Value = 772
$ws4.Cells.Item(20, 6).Value = 1741
$ws4.Cells.Item(21, 6).Value = 2863
$ws4.Cells.Item(22, 6).Value = 833
$ws4.Cells.Item(25, 6).Value = 2199
$ws4.Cells.Item(26, 6).Value = 3040
$ws4.Cells.Item(27, 6).Value = 590
$ws4.Cells.Item(28, 6).Value = 13
$ws4.Cells.Item(30, 6).Value = 12
$ws4.Cells.Item(31, 6).Value = 13
$ws4.Cells.Item(35, 6).Value = 78
$ws4.Cells.Item(36, 6).Value = 734
$ws4.Cells.Item(38, 6).Value = 131
$ws4.Cells.Item(41, 6).Value = 1072
$ws4.Cells.Item(42, 6).Value = 1769
$ws4.Cells.Item(43, 6).Value = 389
$ws4.Cells.Item(44, 6).Value = 554
$ws4.Cells.Item(45, 6).Value = 187
$ws4.Cells.Item(46, 6).Value = 132
$ws4.Cells.Item(47, 6).Value = 173
$ws4.Cells.Item(48, 6).Value = 43
